{"js": "// Replace the date line and each of the division-problem answers in the\n// table with their updated values, per the commit's xml diff.\nconst replacements = [\n  [\"2024-05-10 Friday\", \"2024-05-11 Saturday\"],\n  [\"476\u00f74=119, 0\", \"733\u00f78=91, 5\"],\n  [\"231\u00f78=28, 7\", \"911\u00f76=151, 5\"],\n  [\"464\u00f73=154, 2\", \"523\u00f78=65, 3\"],\n  [\"750\u00f75=150, 0\", \"228\u00f77=32, 4\"],\n  [\"314\u00f74=78, 2\", \"957\u00f78=119, 5\"],\n  [\"513\u00f78=64, 1\", \"352\u00f78=44, 0\"],\n  [\"481\u00f76=80, 1\", \"536\u00f72=268, 0\"],\n  [\"428\u00f78=53, 4\", \"229\u00f78=28, 5\"],\n  [\"930\u00f77=132, 6\", \"913\u00f72=456, 1\"],\n  [\"781\u00f78=97, 5\", \"424\u00f76=70, 4\"],\n  [\"653\u00f75=130, 3\", \"489\u00f79=54, 3\"],\n  [\"476\u00f78=59, 4\", \"287\u00f75=57, 2\"],\n  [\"938\u00f72=469, 0\", \"694\u00f78=86, 6\"],\n  [\"360\u00f72=180, 0\", \"182\u00f77=26, 0\"],\n  [\"892\u00f79=99, 1\", \"533\u00f75=106, 3\"],\n  [\"528\u00f72=264, 0\", \"534\u00f72=267, 0\"],\n  [\"551\u00f72=275, 1\", \"391\u00f75=78, 1\"],\n  [\"952\u00f74=238, 0\", \"909\u00f75=181, 4\"],\n  [\"152\u00f72=76, 0\", \"178\u00f77=25, 3\"],\n  [\"382\u00f74=95, 2\", \"699\u00f74=174, 3\"],\n  [\"505\u00f79=56, 1\", \"669\u00f74=167, 1\"],\n  [\"785\u00f73=261, 2\", \"925\u00f77=132, 1\"],\n  [\"765\u00f73=255, 0\", \"978\u00f74=244, 2\"],\n  [\"199\u00f79=22, 1\", \"817\u00f79=90, 7\"],\n  [\"684\u00f75=136, 4\", \"192\u00f78=24, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each of the division-problem answers in the\n# table with their updated values, per the commit's xml diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-05-10 Friday\", \"2024-05-11 Saturday\"),\n  @(\"476\u00f74=119, 0\", \"733\u00f78=91, 5\"),\n  @(\"231\u00f78=28, 7\", \"911\u00f76=151, 5\"),\n  @(\"464\u00f73=154, 2\", \"523\u00f78=65, 3\"),\n  @(\"750\u00f75=150, 0\", \"228\u00f77=32, 4\"),\n  @(\"314\u00f74=78, 2\", \"957\u00f78=119, 5\"),\n  @(\"513\u00f78=64, 1\", \"352\u00f78=44, 0\"),\n  @(\"481\u00f76=80, 1\", \"536\u00f72=268, 0\"),\n  @(\"428\u00f78=53, 4\", \"229\u00f78=28, 5\"),\n  @(\"930\u00f77=132, 6\", \"913\u00f72=456, 1\"),\n  @(\"781\u00f78=97, 5\", \"424\u00f76=70, 4\"),\n  @(\"653\u00f75=130, 3\", \"489\u00f79=54, 3\"),\n  @(\"476\u00f78=59, 4\", \"287\u00f75=57, 2\"),\n  @(\"938\u00f72=469, 0\", \"694\u00f78=86, 6\"),\n  @(\"360\u00f72=180, 0\", \"182\u00f77=26, 0\"),\n  @(\"892\u00f79=99, 1\", \"533\u00f75=106, 3\"),\n  @(\"528\u00f72=264, 0\", \"534\u00f72=267, 0\"),\n  @(\"551\u00f72=275, 1\", \"391\u00f75=78, 1\"),\n  @(\"952\u00f74=238, 0\", \"909\u00f75=181, 4\"),\n  @(\"152\u00f72=76, 0\", \"178\u00f77=25, 3\"),\n  @(\"382\u00f74=95, 2\", \"699\u00f74=174, 3\"),\n  @(\"505\u00f79=56, 1\", \"669\u00f74=167, 1\"),\n  @(\"785\u00f73=261, 2\", \"925\u00f77=132, 1\"),\n  @(\"765\u00f73=255, 0\", \"978\u00f74=244, 2\"),\n  @(\"199\u00f79=22, 1\", \"817\u00f79=90, 7\"),\n  @(\"684\u00f75=136, 4\", \"192\u00f78=24, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
